# Edit script: adds the new "Problema con los idiomas" section at the end
# of the document body, right after the paragraph containing the image
# "(Imagen de la función con la solución definitiva)".
#
# The rest of the unified diff supplied with this task consists solely of
# Word's own proofing-tool (<w:proofErr>) markup being injected as a side
# effect of the document being re-opened/re-saved with the spell checker
# enabled; the underlying visible text in every other paragraph is
# byte-for-byte identical before and after. Those cosmetic proofErr
# splits cannot be produced deliberately through the Find/Replace/Range
# object model (they are an artifact of Word's spell-check engine at
# save time), so this script focuses on the one substantive content
# change: the new "Problema con los idiomas" paragraphs, matching the
# commit message "Agregado el texto sobre el problema del idioma".

$d = $word.ActiveDocument

# --- Paragraph 1: blank spacer paragraph (style/spacing inherited from
#     the previous paragraph, i.e. "Normal (Web)" with spacing-after 0) ---
$endRange = $d.Content
$endRange.Collapse(0)   # wdCollapseEnd
$endRange.InsertParagraphAfter()

# --- Paragraph 2: section heading "Problema con los idiomas" ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter("Problema con los idiomas")
$headingRange = $d.Paragraphs.Last.Range
$headingRange.Font.Name = "Arial"
$headingRange.Font.Size = 16
$headingRange.Font.Underline = 1

# --- Paragraph 3: body text describing the language bug ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter("A la hora de probar el programa, todo estaba funcionando, entonces probando los idiomas descubrimos que había un error de compilación con el idioma inglés. Al revisarlo mejor descubrimos que señalaba el submenú “métricas” donde en un switch un case había sigo declarado dos veces, lo cual no tenía mucho sentido que el error aparezca solamente cuando cambiamos a ese idioma ya que lo único que hay son definiciones en el idioma, hasta que nos dimos cuenta en el typedef de las variables donde poníamos cual es la letra que hay que ingresar, había una letra que se estaba repitiendo, por lo que para dos casos distintos estabas usando la misma letras, esto se pudo arreglar cambiándolo y ya no se notaba ese problema.")
$bodyRange = $d.Paragraphs.Last.Range
$bodyRange.Font.Name = "Arial"
$bodyRange.Font.Size = 11
$bodyRange.Font.Underline = 0
